$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("CategoricalVariables")

# New category values for iode_quality_flag attribute
$ws.Range("A4").Value = "iode_quality_flag"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "good"

$ws.Range("A5").Value = "iode_quality_flag"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "quality not evaluated, not available or unknown"

$ws.Range("A6").Value = "iode_quality_flag"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "questionable/suspect"

$ws.Range("A7").Value = "iode_quality_flag"
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = "bad"

$ws.Range("A8").Value = "iode_quality_flag"
$ws.Range("B8").Value = 9
$ws.Range("C8").Value = "missing data"

# Make this sheet the active/selected one, then select the new range
$ws.Activate() | Out-Null
$ws.Range("A4:C8").Select() | Out-Null
